$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (name unchanged)
$ws.Range("B3").Value = 0.7560608764112949
$ws.Range("C3").Value = 0.7343812131865159
$ws.Range("D3").Value = 0.652928209961623

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.6515296752250829
$ws.Range("C4").Value = 0.5827910642453481
$ws.Range("D4").Value = 0.606979269069012

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.6652068857830608
$ws.Range("C5").Value = 0.6139830513575016
$ws.Range("D5").Value = 0.271555918505012
